# Add a "Sector" column (F) populated per-ticker, matching the commit
# "Add sector for each ticker".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ticker (column B) -> Sector lookup table.
$sectorMap = @{
    "NVDA"  = "Technology"
    "AAPL"  = "Technology"
    "MSFT"  = "Technology"
    "AMZN"  = "Consumer Cyclical"
    "META"  = "Communication Services"
    "GOOGL" = "Communication Services"
    "TSLA"  = "Consumer Cyclical"
    "GOOG"  = "Communication Services"
    "WMT"   = "Consumer Defensive"
    "AVGO"  = "Technology"
    "LMT"   = "Industrials"
    "TXN"   = "Technology"
    "CVX"   = "Energy"
    "COP"   = "Energy"
    "VZ"    = "Communication Services"
    "BMY"   = "Healthcare"
    "MRK"   = "Healthcare"
    "MO"    = "Consumer Defensive"
    "KO"    = "Consumer Defensive"
    "HD"    = "Consumer Cyclical"
    "LLY"   = "Healthcare"
    "BRK-B" = "Financial Services"
    "MU"    = "Technology"
    "AMD"   = "Technology"
    "PLTR"  = "Technology"
    "CSCO"  = "Technology"
    "LRCX"  = "Technology"
    "IBM"   = "Technology"
}

# Header cell: copy formatting from the existing "holding_value" header (E1)
# so the new column matches the bold/bordered header style, then set text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Sector"

# Fill in the Sector value for every data row based on the Symbol (col B).
$lastRow = 71
for ($row = 2; $row -le $lastRow; $row++) {
    $ticker = $ws.Cells.Item($row, 2).Value()
    if ($sectorMap.ContainsKey($ticker)) {
        $ws.Cells.Item($row, 6).Value = $sectorMap[$ticker]
    }
}
